$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.642273426055908
$ws.Range("B1").Value = 2.871334552764893
$ws.Range("C1").Value = 3.47297739982605
$ws.Range("D1").Value = 3.705272674560547
$ws.Range("E1").Value = 3.30239725112915
